# Generate Report for Handback
# This script updates the localization-status workbook to reflect a
# completed handback: the overall status text changes, per-language
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns are populated, and a hyperlink is added to the new
# "Latest Target File" cell pointing at the source markdown file.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdFileName = "d86cc8b9-4bff-4efe-9b6f-dd23c4759863.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5dbca17c3d468dc6f4e28dfc1e24f9b08ab383e/e2e/d86cc8b9-4bff-4efe-9b6f-dd23c4759863.md"
$zhTargetFile = "d86cc8b9-4bff-4efe-9b6f-dd23c4759863.23e881281bce6773eecbbf5ac61505876f7789d3.zh-cn.xlf"
$deTargetFile = "d86cc8b9-4bff-4efe-9b6f-dd23c4759863.23e881281bce6773eecbbf5ac61505876f7789d3.de-de.xlf"
$zhHandbackDate = "2016-09-03 21:01:48"
$deHandbackDate = "2016-09-03 21:01:55"

# Update the overall status everywhere it is surfaced: the Overview
# sheet's per-language status columns, and each language sheet's own
# Status cell.
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$zh.Range("C2").Value = $newStatus
$de.Range("C2").Value = $newStatus

# Widen the status columns now that the text is longer. (ColumnWidth is
# expressed in characters and gets rounded to the nearest pixel by
# Excel, so the closest representable value is used.)
$ov.Range("E1").EntireColumn.ColumnWidth = 29.1666666666667
$ov.Range("F1").EntireColumn.ColumnWidth = 29.1666666666667
$zh.Range("C1").EntireColumn.ColumnWidth = 29.1666666666667
$de.Range("C1").EntireColumn.ColumnWidth = 29.1666666666667

# Populate "Latest Target File" (I) and "Latest Handback File" (J) for
# each language, and widen those columns to fit the longer file names.
$zh.Range("I2").Value = $mdFileName
$zh.Range("J2").Value = $zhTargetFile
$zh.Range("I1").EntireColumn.ColumnWidth = 39.1666666666667
$zh.Range("J1").EntireColumn.ColumnWidth = 39.1666666666667

$de.Range("I2").Value = $mdFileName
$de.Range("J2").Value = $deTargetFile
$de.Range("I1").EntireColumn.ColumnWidth = 39.1666666666667
$de.Range("J1").EntireColumn.ColumnWidth = 39.1666666666667

# Record the handback timestamps.
$zh.Range("K2").Value = $zhHandbackDate
$de.Range("K2").Value = $deHandbackDate

# Link the new "Latest Target File" cells back to the source markdown
# file, matching the existing hyperlink style used elsewhere.
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)

Write-Host "Handback report generated"
